$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.371.33"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "2.013.39"
$ws.Range("E3").Value = "  -1.70%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.79%  "
$ws.Range("E6").Value = "  -3.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "61.86"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.35%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.21%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  -2.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.907"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").Value = "2.308.41"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +16.23%  "
$ws.Range("E17").Value = "  +1.59%  "
$ws.Range("D18").Value = "1.991.84"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").Value = "36.350.15"
$ws.Range("E19").Value = "  +0.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "0.0₃0866"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "234.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.90%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  +2.43%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.57"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.76%  "
$ws.Range("E29").Value = "  -1.75%  "
$ws.Range("E30").Value = "  -0.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.13%  "
$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.111"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +23.73%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.33%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.58"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0610"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.73%  "
$ws.Range("E36").Value = "  +9.89%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("E38").Value = "  -1.25%  "
$ws.Range("E39").Value = "  +16.94%  "
$ws.Range("E40").Value = "  +15.94%  "
$ws.Range("E41").Value = "  +1.43%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.62%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0217"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.47%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.11"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.37%  "
$ws.Range("B46").Value = "ARBITRUM"
$ws.Range("C46").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.13"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.75%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.33%  "
$ws.Range("D48").Value = "1.454.89"
$ws.Range("E48").Value = "  +4.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.17"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.81%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "47.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.32%  "
